# Add a new "WDI Climate Change" column (column W) to the dataset sheet,
# mirroring the formatting/pattern already used by the other year-coverage
# columns (e.g. column V / "OECD Air Emissions").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column W (rows 3-31) needs the same direct cell formatting (Times New
# Roman 12pt, matching the rest of the table, style index 1) as column A.
# Copy that formatting down first so the new cells line up with the rest
# of the sheet instead of using Excel's bare default style.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("W3:W31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Header (row 3)
$ws.Range("W3").Value = "WDI Climate Change"

# Data rows 4-30: checkmark, same as the rest of the columns for those years
for ($r = 4; $r -le 30; $r++) {
    $ws.Cells.Item($r, 23).Value = "✓"
}

# Row 31 ("Years Vary by Country") is left blank for this dataset, but it
# still carries the formatted (empty) cell like its neighbours.
$ws.Range("W31").Value = $null

# Reflect the new column in the saved view/selection state.
$ws.Range("W31").Select()
